$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply an AutoFilter on the used range (A1:I170), filtering column E
# ("level_3", the 5th column in the range) down to the single discrete
# value "Emprunt" - this is what drives Excel to hide every row whose
# level_3 is not "Emprunt".
$rng = $ws.Range("A1:I170")
$rng.AutoFilter(5, @("Emprunt"), 7)

# Register the (hidden) built-in sheet-scoped defined name Excel creates
# whenever a range has an AutoFilter applied to it.
$filterDbName = $ws.Names.Add("_xlnm._FilterDatabase", "='mappings_2026-01-19'!`$A`$1:`$I`$170")
$filterDbName.Visible = $false

# Match the author's resulting selection: columns B:D, active cell B1.
$ws.Range("B1:D1048576").Select()
